# Updates numeric leve-profit figures (currentAveragePrice / Leve price / profit columns)
# across the per-class worksheets, per the scheduled Odin_Profits data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()  # was -1380
$ws.Range("H42").Value = 698.7
$ws.Range("I42").Value = 1425.75
$ws.Range("J42").Value = 214
$ws.Range("K42").Value = 4277.25
$ws.Range("L42").Value = 642
$ws.Range("M42").Value = -4047.25
$ws.Range("N42").Value = -1102
$ws.Range("H76").Value = 250008750
$ws.Range("J76").Value = 15000
$ws.Range("L76").Value = 15000
$ws.Range("N76").Value = -15630
$ws.Range("H79").Value = 250008750
$ws.Range("J79").Value = 15000
$ws.Range("L79").Value = 15000
$ws.Range("N79").Value = -17184
$ws.Range("H82").Value = 4532.5
$ws.Range("I82").Value = 4532.5
$ws.Range("K82").Value = 13597.5
$ws.Range("M82").Value = -13191.5
$ws.Range("H85").Value = 4532.5
$ws.Range("I85").Value = 4532.5
$ws.Range("K85").Value = 13597.5
$ws.Range("M85").Value = -12193.5
$ws.Range("H92").Value = 1261.25
$ws.Range("I92").Value = 1180.8823
$ws.Range("K92").Value = 1180.8823
$ws.Range("M92").Value = 67.11770000000001
$ws.Range("H100").Value = 1299.875
$ws.Range("I100").Value = 1128.4286
$ws.Range("K100").Value = 1128.4286
$ws.Range("M100").Value = -587.4286
$ws.Range("H101").Value = 3196.125
$ws.Range("I101").Value = 2224.2856
$ws.Range("K101").Value = 6672.8568
$ws.Range("M101").Value = -5050.8568
$ws.Range("H113").Value = 900000
$ws.Range("I113").Value = 900000
$ws.Range("J113").Value = 900000
$ws.Range("K113").Value = 900000
$ws.Range("L113").Value = 900000
$ws.Range("M113").Value = -896746
$ws.Range("N113").Value = -906508
$ws.Range("H115").Value = 2064.7
$ws.Range("I115").Value = 1307.1428
$ws.Range("J115").Value = 3832.3333
$ws.Range("K115").Value = 3921.4284
$ws.Range("L115").Value = 11496.9999
$ws.Range("M115").Value = -2354.4284
$ws.Range("N115").Value = -14630.9999
$ws.Range("H135").Value = 4939.3237
$ws.Range("I135").Value = 1668.5834
$ws.Range("K135").Value = 15017.2506
$ws.Range("M135").Value = -12482.2506
$ws.Range("H137").Value = 12857.214
$ws.Range("I137").Value = 13487.5
$ws.Range("J137").Value = 12605.1
$ws.Range("K137").Value = 40462.5
$ws.Range("L137").Value = 37815.3
$ws.Range("M137").Value = -37912.5
$ws.Range("N137").Value = -42915.3
$ws.Range("H138").Value = 3144.1694
$ws.Range("I138").Value = 821.6667
$ws.Range("J138").Value = 5103.7812
$ws.Range("K138").Value = 2465.0001
$ws.Range("L138").Value = 15311.3436
$ws.Range("M138").Value = 2674.9999
$ws.Range("N138").Value = -25591.3436

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5114.9443
$ws.Range("I2").Value = 2644.3076
$ws.Range("J2").Value = 11538.6
$ws.Range("K2").Value = 2644.3076
$ws.Range("L2").Value = 11538.6
$ws.Range("M2").Value = -2531.3076
$ws.Range("N2").Value = -11764.6
$ws.Range("H32").Value = 2714011.5
$ws.Range("I32").Value = 1614.7693
$ws.Range("K32").Value = 1614.7693
$ws.Range("M32").Value = -1327.7693
$ws.Range("H97").Value = 1991
$ws.Range("I97").Value = 1624.5
$ws.Range("J97").Value = 2409.8572
$ws.Range("K97").Value = 1624.5
$ws.Range("L97").Value = 2409.8572
$ws.Range("M97").Value = -1128.5
$ws.Range("N97").Value = -3401.8572
$ws.Range("H116").Value = 5114.9443
$ws.Range("I116").Value = 2644.3076
$ws.Range("J116").Value = 11538.6
$ws.Range("K116").Value = 2644.3076
$ws.Range("L116").Value = 11538.6
$ws.Range("M116").Value = -350.3076000000001
$ws.Range("N116").Value = -16126.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5114.9443
$ws.Range("I3").Value = 2644.3076
$ws.Range("J3").Value = 11538.6
$ws.Range("K3").Value = 2644.3076
$ws.Range("L3").Value = 11538.6
$ws.Range("M3").Value = -2530.3076
$ws.Range("N3").Value = -11766.6
$ws.Range("H94").Value = 6262.5137
$ws.Range("I94").Value = 3296.9167
$ws.Range("K94").Value = 3296.9167
$ws.Range("M94").Value = -2845.9167
$ws.Range("H99").Value = 7877.887
$ws.Range("J99").Value = 7494.926
$ws.Range("L99").Value = 7494.926
$ws.Range("N99").Value = -10490.926

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3860.25
$ws.Range("I31").Value = 1830.9048
$ws.Range("K31").Value = 1830.9048
$ws.Range("M31").Value = -1535.9048
$ws.Range("H34").Value = 3860.25
$ws.Range("I34").Value = 1830.9048
$ws.Range("K34").Value = 1830.9048
$ws.Range("M34").Value = -1628.9048
$ws.Range("H56").Value = 31565
$ws.Range("J56").Value = 33130
$ws.Range("L56").Value = 33130
$ws.Range("N56").Value = -34820

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 4000061.5
$ws.Range("J12").Value = 79.333336
$ws.Range("L12").Value = 238.000008
$ws.Range("N12").Value = -584.000008
$ws.Range("H25").Value = 1730
$ws.Range("I25").Value = 1950
$ws.Range("J25").Value = 1400
$ws.Range("K25").Value = 5850
$ws.Range("L25").Value = 4200
$ws.Range("M25").Value = -5681
$ws.Range("N25").Value = -4538
$ws.Range("H26").Value = 397.93332
$ws.Range("I26").Value = 432.5
$ws.Range("J26").Value = 259.66666
$ws.Range("K26").Value = 1297.5
$ws.Range("L26").Value = 778.9999799999999
$ws.Range("M26").Value = -1009.5
$ws.Range("N26").Value = -1354.99998
$ws.Range("H30").Value = 1730
$ws.Range("I30").Value = 1950
$ws.Range("J30").Value = 1400
$ws.Range("K30").Value = 5850
$ws.Range("L30").Value = 4200
$ws.Range("M30").Value = -5748
$ws.Range("N30").Value = -4404
$ws.Range("H37").Value = 135997
$ws.Range("J37").Value = 135997
$ws.Range("L37").Value = 407991
$ws.Range("N37").Value = -408215
$ws.Range("H52").Value = 2315518.2
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 2315518.2
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 6946554.600000001
$ws.Range("M52").ClearContents()  # was -334
$ws.Range("N52").Value = -6947086.600000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 39996.668
$ws.Range("J53").Value = 39996.668
$ws.Range("L53").Value = 39996.668
$ws.Range("N53").Value = -41258.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 18195.8
$ws.Range("J39").Value = 18195.8
$ws.Range("L39").Value = 18195.8
$ws.Range("N39").Value = -19115.8
$ws.Range("H40").Value = 4601.718
$ws.Range("I40").Value = 3597.1035
$ws.Range("K40").Value = 3597.1035
$ws.Range("M40").Value = -3461.1035
$ws.Range("H61").Value = 5376.2964
$ws.Range("I61").Value = 4589.1333
$ws.Range("J61").Value = 6360.25
$ws.Range("K61").Value = 4589.1333
$ws.Range("L61").Value = 6360.25
$ws.Range("M61").Value = -4387.1333
$ws.Range("N61").Value = -6764.25
$ws.Range("H68").Value = 2369.926
$ws.Range("I68").Value = 1941.8422
$ws.Range("J68").Value = 3386.625
$ws.Range("K68").Value = 1941.8422
$ws.Range("L68").Value = 3386.625
$ws.Range("M68").Value = -1192.8422
$ws.Range("N68").Value = -4884.625
$ws.Range("H71").Value = 2369.926
$ws.Range("I71").Value = 1941.8422
$ws.Range("J71").Value = 3386.625
$ws.Range("K71").Value = 9709.210999999999
$ws.Range("L71").Value = 16933.125
$ws.Range("M71").Value = -5965.210999999999
$ws.Range("N71").Value = -24421.125
$ws.Range("H106").Value = 27361.5
$ws.Range("J106").Value = 27361.5
$ws.Range("L106").Value = 27361.5
$ws.Range("N106").Value = -29885.5
$ws.Range("H113").Value = 5376.2964
$ws.Range("I113").Value = 4589.1333
$ws.Range("J113").Value = 6360.25
$ws.Range("K113").Value = 4589.1333
$ws.Range("L113").Value = 6360.25
$ws.Range("M113").Value = -2419.1333
$ws.Range("N113").Value = -10700.25
$ws.Range("H136").Value = 107150480
$ws.Range("I136").Value = 50007676
$ws.Range("K136").Value = 150023028
$ws.Range("M136").Value = -150020478

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2080.2727
$ws.Range("I100").Value = 874.75
$ws.Range("K100").Value = 1749.5
$ws.Range("M100").Value = -1208.5
$ws.Range("H107").Value = 520.3043
$ws.Range("I107").Value = 526.7619
$ws.Range("J107").Value = 452.5
$ws.Range("K107").Value = 1580.2857
$ws.Range("L107").Value = 1357.5
$ws.Range("M107").Value = 339.7143000000001
$ws.Range("N107").Value = -5197.5
$ws.Range("H113").Value = 5953209.5
$ws.Range("I113").Value = 9805009
$ws.Range("K113").Value = 29415027
$ws.Range("M113").Value = -29412857
$ws.Range("H132").Value = 6990.42
$ws.Range("I132").Value = 5413.0835
$ws.Range("J132").Value = 11046.429
$ws.Range("K132").Value = 16239.2505
$ws.Range("L132").Value = 33139.287
$ws.Range("M132").Value = -13709.2505
$ws.Range("N132").Value = -38199.287

